{"js": "// Rewrite the hypothesis/graph paragraph and tighten the axis description,\n// then relocate the \"_GoBack\" bookmark to sit after \"and the X axis\" in the\n// third paragraph (it previously sat at the very end of the document).\n\nconst body = context.document.body;\n\nasync function replaceOnce(find, replace) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n  results.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n\n// Paragraph 2 (\"Our initial hypothesis...\") wording tweaks.\nawait replaceOnce(\"more grim, moving\", \"grimmer as Dante moves\");\nawait replaceOnce(\"this trend is only a proof\", \"this graph is a proof\");\n\n// Paragraph 3 (\"The graph displays...\") \u2014 drop the parenthetical asides.\nawait replaceOnce(\" (vertical axis)\", \"\");\nawait replaceOnce(\" (horizontal axis)\", \"\");\n\n// Move the \"_GoBack\" bookmark from the end of the document to right after\n// \"and the X axis\" in paragraph 3.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst anchor = body.search(\"and the X axis\", { matchCase: true });\nanchor.load(\"text\");\nawait context.sync();\nif (anchor.items.length === 0) {\n  throw new Error(\"Bookmark anchor text not found\");\n}\nconst anchorEnd = anchor.items[0].getRange(\"End\");\nanchorEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Rewrite the hypothesis/graph paragraph and tighten the axis description,\n# then relocate the \"_GoBack\" bookmark to sit after \"and the X axis\" in the\n# third paragraph (it previously sat at the very end of the document).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.Text = $findText\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n    $range.Text = $replaceText\n}\n\n# Paragraph 2 (\"Our initial hypothesis...\") wording tweaks.\nReplace-Text \"more grim, moving\" \"grimmer as Dante moves\"\nReplace-Text \"this trend is only a proof\" \"this graph is a proof\"\n\n# Paragraph 3 (\"The graph displays...\") - drop the parenthetical asides.\nReplace-Text \" (vertical axis)\" \"\"\nReplace-Text \" (horizontal axis)\" \"\"\n\n# Move the \"_GoBack\" bookmark from the end of the document to right after\n# \"and the X axis\" in paragraph 3.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$anchorFind = $d.Content.Find\n$anchorFind.Text = \"and the X axis\"\n$anchorFound = $anchorFind.Execute()\nif (-not $anchorFound) {\n    throw \"Bookmark anchor text not found\"\n}\n$anchorEnd = $anchorFind.Parent.Duplicate\n$anchorEnd.Collapse(0)  # wdCollapseEnd\n\n$d.Bookmarks.Add(\"_GoBack\", $anchorEnd)\n"}
